$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell captions that changed wording (R&R edits).
$ws.Range("A5").Value = "Mandatory structured (ATE)"
$ws.Range("A7").Value = "Choice  (ITT)"
$ws.Range("B3").Value = 'Choose structure in $t+1$'
$ws.Range("C3").Value = 'Ever choose structure in $t+1$'

# Scroll the view back to the top-left (A1) instead of A5.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
